# Comparison.xlsx — update the "old" sample column (B) so that each value
# matches its recomputed ("new") counterpart already present in column C,
# and refresh the active sheet's view state (zoom + selection) to match
# where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B corrections (B == C after the edit) ---
$ws.Range("B2").Value  = 92
$ws.Range("B3").Value  = 4
$ws.Range("B4").Value  = 240
$ws.Range("B6").Value  = 118
$ws.Range("B7").Value  = 447
$ws.Range("B8").Value  = 53
$ws.Range("B9").Value  = 148
$ws.Range("B10").Value = 541
$ws.Range("B12").Value = 153
$ws.Range("B15").Value = 223

# --- View state: selection moves from B15 to C15, zoom to 160% ---
$ws.Range("C15").Select()
$excel.ActiveWindow.Zoom = 160
